# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose B:AC contents (odds/results data) were swapped between
# each other while keeping the row index (column A) in place.
$pairs = @(
    @(171,172),
    @(181,182),
    @(186,187),
    @(197,198),
    @(233,234),
    @(241,242),
    @(250,251),
    @(263,264),
    @(273,274),
    @(281,282),
    @(317,318)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $v1 = $range1.Value()
    $v2 = $range2.Value()

    $range1.Value = $v2
    $range2.Value = $v1
}

# The last two data rows (319 and 320) were removed entirely.
$ws.Rows.Item(320).Delete()
$ws.Rows.Item(319).Delete()
